$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.07131444737854614
$ws.Range("J2").Value = 0.07131444737854616
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 13.57958433333333
$ws.Range("N2").Value = 40.738753
$ws.Range("O2").Value = 0.2289698008477291
$ws.Range("P2").Value = 0.2289698008477291
$ws.Range("Q2").Value = 0.368319065873
$ws.Range("R2").Value = 3.314871592857
$ws.Range("S2").Value = 0.01632885481383156
$ws.Range("T2").Value = 0.01632885481383156

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.07131444737854614
$ws.Range("J3").Value = 0.07131444737854616
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.768727
$ws.Range("N3").Value = 59.306181
$ws.Range("O3").Value = 0.3333269541315948
$ws.Range("P3").Value = 0.3333269541315948
$ws.Range("Q3").Value = 0.5361871824209999
$ws.Range("R3").Value = 4.825684641789
$ws.Range("S3").Value = 0.02377102753026868
$ws.Range("T3").Value = 0.02377102753026869

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.07131444737854614
$ws.Range("J4").Value = 0.07131444737854616
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 25.95900466666667
$ws.Range("N4").Value = 77.877014
$ws.Range("O4").Value = 0.4377032450206762
$ws.Range("P4").Value = 0.4377032450206762
$ws.Range("Q4").Value = 0.704086083574
$ws.Range("R4").Value = 6.336774752166
$ws.Range("S4").Value = 0.0312145650344459
$ws.Range("T4").Value = 0.03121456503444591

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3532066666666667
$ws.Range("H5").Value = 1.05962
$ws.Range("I5").Value = 0.9286855526214538
$ws.Range("J5").Value = 0.9286855526214538
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 13.57958433333333
$ws.Range("N5").Value = 40.738753
$ws.Range("O5").Value = 0.2289698008477291
$ws.Range("P5").Value = 0.2289698008477291
$ws.Range("Q5").Value = 4.796399717095556
$ws.Range("R5").Value = 43.16759745386
$ws.Range("S5").Value = 0.2126409460338975
$ws.Range("T5").Value = 0.2126409460338975

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3532066666666667
$ws.Range("H6").Value = 1.05962
$ws.Range("I6").Value = 0.9286855526214538
$ws.Range("J6").Value = 0.9286855526214538
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.768727
$ws.Range("N6").Value = 59.306181
$ws.Range("O6").Value = 0.3333269541315948
$ws.Range("P6").Value = 0.3333269541315948
$ws.Range("Q6").Value = 6.982446167913333
$ws.Range("R6").Value = 62.84201551122
$ws.Range("S6").Value = 0.3095559266013261
$ws.Range("T6").Value = 0.3095559266013261

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3532066666666667
$ws.Range("H7").Value = 1.05962
$ws.Range("I7").Value = 0.9286855526214538
$ws.Range("J7").Value = 0.9286855526214538
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.95900466666667
$ws.Range("N7").Value = 77.877014
$ws.Range("O7").Value = 0.4377032450206762
$ws.Range("P7").Value = 0.4377032450206762
$ws.Range("Q7").Value = 9.168893508297778
$ws.Range("R7").Value = 82.52004157468001
$ws.Range("S7").Value = 0.4064886799862303
$ws.Range("T7").Value = 0.4064886799862303
